$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New mini lookup table (rows 67-73), a condensed copy of the
#     Name/Type1/Total-stats columns used later by INDEX/MATCH examples ---
$ws.Range("A67").Value = "    Name"
$ws.Range("A67").WrapText = $true
$ws.Range("B67").Value = "Type 1"
$ws.Range("C67").Value = "Total stats"

$ws.Range("A68").Value = "    Mankey"
$ws.Range("B68").Value = "Fighting"
$ws.Range("C68").Value = 305

$ws.Range("A69").Value = "    Poliwrath"
$ws.Range("B69").Value = "Water"
$ws.Range("C69").Value = 510

$ws.Range("A70").Value = "    Victreebel"
$ws.Range("B70").Value = "Grass"
$ws.Range("C70").Value = 490

$ws.Range("A71").Value = "    Tentacool"
$ws.Range("B71").Value = "Water"
$ws.Range("C71").Value = 335

$ws.Range("A72").Value = "    Magneton"
$ws.Range("B72").Value = "Electric"
$ws.Range("C72").Value = 465

$ws.Range("A73").Value = "    Dewgong"
$ws.Range("B73").Value = "Water"
$ws.Range("C73").Value = 475

# --- Left Function example ---
$ws.Range("A78").Value = "Left Function"

$ws.Range("A79").Formula = "=LEFT(B68)"
$ws.Range("A80:A85").Formula = "=LEFT(B69)"

# --- Index function example ---
$ws.Range("A86").Value = "index function"

$ws.Range("A87").Formula = "=INDEX(A67:C73,4,2)"
$ws.Range("A88").Formula = "=INDEX(A67:C73,3,2)"

# --- Match function heading ---
$ws.Range("A91").Value = "match function"

# --- Update the view's scroll/selection to reflect the new content ---
$excel.ActiveWindow.ScrollRow = 65
$ws.Range("A92").Select()
